$d = $word.ActiveDocument

# 1) Professional summary paragraph: simple text swap.
$d.Content.Find.Execute(
    "affecting all Black and Asian-American voters, developed geospatial ML algorithms",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "affecting 50M voters, developed geospatial ML algorithms", 2)

# 2) Work-experience bullet: needs "50M" as its own bold/colored run, with
#    " voters," preserved as plain text afterwards. Locate the specific
#    paragraph so we only touch this occurrence.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Discovered systematic race coding errors*") {
        $target = $p.Range
    }
}
$target.Find.Execute("all Black and Asian-American", $false, $false, $false, $false, $false, $true, 1, $false, "", 0, $false)
$target.Text = "50M"
$target.Bold = 1
$target.Font.Color = 5258796  # wdColor RGB(0x2C,0x3E,0x50) -> matches w:color val="2C3E50"

# 3) Key Projects "Impact:" line: simple text swap, adding "nationwide".
$d.Content.Find.Execute(
    "affecting all Black and Asian-American voters, improved electoral prediction accuracy by 22%",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "affecting 50M voters nationwide, improved electoral prediction accuracy by 22%", 2)
